$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5546926666666666
$ws.Range("H2").Value = 1.664078
$ws.Range("I2").Value = 0.01048496710804503
$ws.Range("J2").Value = 0.01048496710804503
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1747913333333333
$ws.Range("N2").Value = 0.524374
$ws.Range("O2").Value = 0.08331024341015711
$ws.Range("P2").Value = 0.08331024341015709
$ws.Range("Q2").Value = 0.09695547079688888
$ws.Range("R2").Value = 0.8725992371719999
$ws.Range("S2").Value = 0.0008735051619187223
$ws.Range("T2").Value = 0.0008735051619187223
$ws.Range("G3").Value = 0.5546926666666666
$ws.Range("H3").Value = 1.664078
$ws.Range("I3").Value = 0.01048496710804503
$ws.Range("J3").Value = 0.01048496710804503
$ws.Range("O3").Value = 0.4977364990677179
$ws.Range("P3").Value = 0.4977364990677178
$ws.Range("Q3").Value = 0.5792598199757778
$ws.Range("R3").Value = 5.213338379782001
$ws.Range("S3").Value = 0.005218750821198507
$ws.Range("T3").Value = 0.005218750821198507
$ws.Range("G4").Value = 0.5546926666666666
$ws.Range("H4").Value = 1.664078
$ws.Range("I4").Value = 0.01048496710804503
$ws.Range("J4").Value = 0.01048496710804503
$ws.Range("M4").Value = 0.8789963333333334
$ws.Range("O4").Value = 0.418953257522125
$ws.Range("P4").Value = 0.418953257522125
$ws.Range("Q4").Value = 0.4875728201268889
$ws.Range("R4").Value = 4.388155381142
$ws.Range("S4").Value = 0.004392711124927799
$ws.Range("T4").Value = 0.004392711124927799
$ws.Range("I5").Value = 0.3689509033301984
$ws.Range("J5").Value = 0.3689509033301984
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1747913333333333
$ws.Range("N5").Value = 0.524374
$ws.Range("O5").Value = 0.08331024341015711
$ws.Range("P5").Value = 0.08331024341015709
$ws.Range("Q5").Value = 3.411723486082222
$ws.Range("R5").Value = 30.70551137474
$ws.Range("S5").Value = 0.03073738956283617
$ws.Range("T5").Value = 0.03073738956283617
$ws.Range("I6").Value = 0.3689509033301984
$ws.Range("J6").Value = 0.3689509033301984
$ws.Range("O6").Value = 0.4977364990677179
$ws.Range("P6").Value = 0.4977364990677178
$ws.Range("S6").Value = 0.183640330951445
$ws.Range("T6").Value = 0.1836403309514449
$ws.Range("I7").Value = 0.3689509033301984
$ws.Range("J7").Value = 0.3689509033301984
$ws.Range("M7").Value = 0.8789963333333334
$ws.Range("O7").Value = 0.418953257522125
$ws.Range("P7").Value = 0.418953257522125
$ws.Range("S7").Value = 0.1545731828159173
$ws.Range("T7").Value = 0.1545731828159172
$ws.Range("I8").Value = 0.6205641295617567
$ws.Range("J8").Value = 0.6205641295617566
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1747913333333333
$ws.Range("N8").Value = 0.524374
$ws.Range("O8").Value = 0.08331024341015711
$ws.Range("P8").Value = 0.08331024341015709
$ws.Range("Q8").Value = 5.738414505388
$ws.Range("R8").Value = 51.645730548492
$ws.Range("S8").Value = 0.05169934868540222
$ws.Range("T8").Value = 0.0516993486854022
$ws.Range("I9").Value = 0.6205641295617567
$ws.Range("J9").Value = 0.6205641295617566
$ws.Range("O9").Value = 0.4977364990677179
$ws.Range("P9").Value = 0.4977364990677178
$ws.Range("S9").Value = 0.3088774172950745
$ws.Range("T9").Value = 0.3088774172950743
$ws.Range("I10").Value = 0.6205641295617567
$ws.Range("J10").Value = 0.6205641295617566
$ws.Range("M10").Value = 0.8789963333333334
$ws.Range("O10").Value = 0.418953257522125
$ws.Range("P10").Value = 0.418953257522125
$ws.Range("Q10").Value = 0.25998736358128
$ws.Range("R10").Value = 0.2599873635812799
$ws.Range("S10").Value = 0.25998736358128
$ws.Range("T10").Value = 0.2599873635812799
